$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.575835466384888
$ws.Range("B1").Value = 7.313015460968018
$ws.Range("C1").Value = 7.071954250335693
$ws.Range("D1").Value = 6.377043724060059
$ws.Range("E1").Value = 3.362801551818848
